# Append 8 new daily data rows (77-84, dated 2024-08-03 .. 2024-08-10) to Sheet1,
# extending the used range from A1:Z76 to A1:Z84.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style (including date number format YYYY-MM-DD HH:MM:SS) from A76 down to A77:A84
$ws.Range("A76").Copy()
$ws.Range("A77:A84").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 77
$ws.Cells.Item(77, 1).Value = 45507
$ws.Cells.Item(77, 2).Value = 771.7209216389
$ws.Cells.Item(77, 3).Value = 201.327799314
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(77, 5).Value = 0
$ws.Cells.Item(77, 6).Value = 0
$ws.Cells.Item(77, 7).Value = 0
$ws.Cells.Item(77, 9).Value = 209.1019791028
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 11).Value = 0.05510130712
$ws.Cells.Item(77, 12).Value = 0
$ws.Cells.Item(77, 13).Value = 0
$ws.Cells.Item(77, 14).Value = 107.43841100064
$ws.Cells.Item(77, 15).Value = 53.807502992
$ws.Cells.Item(77, 16).Value = 0
$ws.Cells.Item(77, 17).Value = 0.0000020688
$ws.Cells.Item(77, 18).Value = 0
$ws.Cells.Item(77, 19).Value = 0
$ws.Cells.Item(77, 20).Value = 0
$ws.Cells.Item(77, 21).Value = 303.5755600274733
$ws.Cells.Item(77, 23).Value = 0
$ws.Cells.Item(77, 24).Value = 0
$ws.Cells.Item(77, 25).Value = 0
$ws.Cells.Item(77, 26).Value = 217.53019326057

# Row 78
$ws.Cells.Item(78, 1).Value = 45508
$ws.Cells.Item(78, 2).Value = 739.46535171
$ws.Cells.Item(78, 3).Value = 186.439898242
$ws.Cells.Item(78, 4).Value = 0
$ws.Cells.Item(78, 5).Value = 0
$ws.Cells.Item(78, 6).Value = 0
$ws.Cells.Item(78, 7).Value = 0
$ws.Cells.Item(78, 9).Value = 202.9600183293
$ws.Cells.Item(78, 10).Value = 0
$ws.Cells.Item(78, 11).Value = 0.05670517489600001
$ws.Cells.Item(78, 12).Value = 0
$ws.Cells.Item(78, 13).Value = 0
$ws.Cells.Item(78, 14).Value = 102.81524270688
$ws.Cells.Item(78, 15).Value = 50.485174163
$ws.Cells.Item(78, 16).Value = 0
$ws.Cells.Item(78, 17).Value = 0.000001944
$ws.Cells.Item(78, 18).Value = 0
$ws.Cells.Item(78, 19).Value = 0
$ws.Cells.Item(78, 20).Value = 0
$ws.Cells.Item(78, 21).Value = 287.5844327609608
$ws.Cells.Item(78, 23).Value = 0
$ws.Cells.Item(78, 24).Value = 0
$ws.Cells.Item(78, 25).Value = 0
$ws.Cells.Item(78, 26).Value = 204.047690607978

# Row 79
$ws.Cells.Item(79, 1).Value = 45509
$ws.Cells.Item(79, 2).Value = 686.8010924091
$ws.Cells.Item(79, 3).Value = 167.7655390965
$ws.Cells.Item(79, 4).Value = 0
$ws.Cells.Item(79, 5).Value = 0
$ws.Cells.Item(79, 6).Value = 0
$ws.Cells.Item(79, 7).Value = 0
$ws.Cells.Item(79, 9).Value = 190.4204875734
$ws.Cells.Item(79, 10).Value = 0
$ws.Cells.Item(79, 11).Value = 0.05566881411
$ws.Cells.Item(79, 12).Value = 0
$ws.Cells.Item(79, 13).Value = 0
$ws.Cells.Item(79, 14).Value = 90.41979148447999
$ws.Cells.Item(79, 15).Value = 47.162845334
$ws.Cells.Item(79, 16).Value = 0
$ws.Cells.Item(79, 17).Value = 0.0000016896
$ws.Cells.Item(79, 18).Value = 0
$ws.Cells.Item(79, 19).Value = 0
$ws.Cells.Item(79, 20).Value = 0
$ws.Cells.Item(79, 21).Value = 243.7047795416505
$ws.Cells.Item(79, 23).Value = 0
$ws.Cells.Item(79, 24).Value = 0
$ws.Cells.Item(79, 25).Value = 0
$ws.Cells.Item(79, 26).Value = 212.370223109578

# Row 80
$ws.Cells.Item(80, 1).Value = 45510
$ws.Cells.Item(80, 2).Value = 712.2699975611
$ws.Cells.Item(80, 3).Value = 170.6596383455
$ws.Cells.Item(80, 4).Value = 0
$ws.Cells.Item(80, 5).Value = 0
$ws.Cells.Item(80, 6).Value = 0
$ws.Cells.Item(80, 7).Value = 0
$ws.Cells.Item(80, 9).Value = 211.87752074
$ws.Cells.Item(80, 10).Value = 0
$ws.Cells.Item(80, 11).Value = 0.057719386246
$ws.Cells.Item(80, 12).Value = 0
$ws.Cells.Item(80, 13).Value = 0
$ws.Cells.Item(80, 14).Value = 98.39308173024001
$ws.Cells.Item(80, 15).Value = 49.25581089600001
$ws.Cells.Item(80, 16).Value = 0
$ws.Cells.Item(80, 17).Value = 0.0000018576
$ws.Cells.Item(80, 18).Value = 0
$ws.Cells.Item(80, 19).Value = 0
$ws.Cells.Item(80, 20).Value = 0
$ws.Cells.Item(80, 21).Value = 271.7212345125803
$ws.Cells.Item(80, 23).Value = 0
$ws.Cells.Item(80, 24).Value = 0
$ws.Cells.Item(80, 25).Value = 0
$ws.Cells.Item(80, 26).Value = 214.076342272406

# Row 81
$ws.Cells.Item(81, 1).Value = 45511
$ws.Cells.Item(81, 2).Value = 700.9817749976
$ws.Cells.Item(81, 3).Value = 162.44120078
$ws.Cells.Item(81, 4).Value = 0
$ws.Cells.Item(81, 5).Value = 0
$ws.Cells.Item(81, 6).Value = 0
$ws.Cells.Item(81, 7).Value = 0
$ws.Cells.Item(81, 9).Value = 212.4261849459
$ws.Cells.Item(81, 10).Value = 0
$ws.Cells.Item(81, 11).Value = 0.059519913548
$ws.Cells.Item(81, 12).Value = 0
$ws.Cells.Item(81, 13).Value = 0
$ws.Cells.Item(81, 14).Value = 91.82684270432
$ws.Cells.Item(81, 15).Value = 48.26012825
$ws.Cells.Item(81, 16).Value = 0
$ws.Cells.Item(81, 17).Value = 0.0000017664
$ws.Cells.Item(81, 18).Value = 0
$ws.Cells.Item(81, 19).Value = 0
$ws.Cells.Item(81, 20).Value = 0
$ws.Cells.Item(81, 21).Value = 260.7193389532198
$ws.Cells.Item(81, 23).Value = 0
$ws.Cells.Item(81, 24).Value = 0
$ws.Cells.Item(81, 25).Value = 0
$ws.Cells.Item(81, 26).Value = 201.446899201228

# Row 82
$ws.Cells.Item(82, 1).Value = 45512
$ws.Cells.Item(82, 2).Value = 784.2824623189
$ws.Cells.Item(82, 3).Value = 185.994758875
$ws.Cells.Item(82, 4).Value = 0
$ws.Cells.Item(82, 5).Value = 0
$ws.Cells.Item(82, 6).Value = 0
$ws.Cells.Item(82, 7).Value = 0
$ws.Cells.Item(82, 9).Value = 239.3729636837
$ws.Cells.Item(82, 10).Value = 0
$ws.Cells.Item(82, 11).Value = 0.06183217123
$ws.Cells.Item(82, 12).Value = 0
$ws.Cells.Item(82, 13).Value = 0
$ws.Cells.Item(82, 14).Value = 110.62102685504
$ws.Cells.Item(82, 15).Value = 52.53749961700001
$ws.Cells.Item(82, 16).Value = 0
$ws.Cells.Item(82, 17).Value = 0.0000021048
$ws.Cells.Item(82, 18).Value = 0
$ws.Cells.Item(82, 19).Value = 0
$ws.Cells.Item(82, 20).Value = 0
$ws.Cells.Item(82, 21).Value = 295.6439609032831
$ws.Cells.Item(82, 23).Value = 0
$ws.Cells.Item(82, 24).Value = 0
$ws.Cells.Item(82, 25).Value = 0
$ws.Cells.Item(82, 26).Value = 215.824074097742

# Row 83
$ws.Cells.Item(83, 1).Value = 45513
$ws.Cells.Item(83, 3).Value = 180.189919653
$ws.Cells.Item(83, 4).Value = 0
$ws.Cells.Item(83, 5).Value = 0
$ws.Cells.Item(83, 7).Value = 0
$ws.Cells.Item(83, 10).Value = 0
$ws.Cells.Item(83, 11).Value = 0.063691986962
$ws.Cells.Item(83, 12).Value = 0
$ws.Cells.Item(83, 13).Value = 0
$ws.Cells.Item(83, 14).Value = 108.34294392768
$ws.Cells.Item(83, 16).Value = 0
$ws.Cells.Item(83, 20).Value = 0
$ws.Cells.Item(83, 24).Value = 0
$ws.Cells.Item(83, 26).Value = 206.731707339744

# Row 84
$ws.Cells.Item(84, 1).Value = 45514
$ws.Cells.Item(84, 3).Value = 180.962326592
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(84, 5).Value = 0
$ws.Cells.Item(84, 7).Value = 0
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 11).Value = 0.065323669189
$ws.Cells.Item(84, 12).Value = 0
$ws.Cells.Item(84, 13).Value = 0
$ws.Cells.Item(84, 14).Value = 106.66788295168
$ws.Cells.Item(84, 16).Value = 0
$ws.Cells.Item(84, 17).Value = 0.0000020256
$ws.Cells.Item(84, 19).Value = 0
$ws.Cells.Item(84, 20).Value = 0
$ws.Cells.Item(84, 24).Value = 0
$ws.Cells.Item(84, 26).Value = 223.56402932423
